$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new header cells (Wins, Losses, Ties) in AD1:AF1, matching the
# formatting (bold, bordered, centered) already used by the other header
# cells on row 1 -- copy the format from an existing header cell (AC1) and
# then set the text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every player row, 2-50.
$lastRow = 50
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 94
    $ws.Cells.Item($r, 31).Value = 67
    $ws.Cells.Item($r, 32).Value = 0
}
